# "inititialisation capture mobile et teste mobile"
# Adds a new "mobile" worksheet (after "bureau") with a mobile test-case
# table, restyles the new header/highlight cells, and updates the view
# state (selection / active tab) on both sheets.

$wb = $excel.ActiveWorkbook
$bureau = $wb.Worksheets.Item("bureau")

# --- "bureau" sheet view: no longer the selected tab, selection moves to C5 ---
$bureau.Range("C5").Select()

# --- create the new "mobile" sheet right after "bureau" ---
$mobile = $wb.Worksheets.Add($null, $bureau)
$mobile.Name = "mobile"

# Column widths (A:D)
$mobile.Columns.Item(1).ColumnWidth = 43.333333
$mobile.Columns.Item(2).ColumnWidth = 51
$mobile.Columns.Item(3).ColumnWidth = 52
$mobile.Columns.Item(4).ColumnWidth = 30

# --- Row 1 : title ---
$mobile.Range("A1").Value = "Test cases"
$mobile.Range("B1").Value = "mobile"

# --- Row 4 : Démarrage de l'app ---
$mobile.Range("A4").Value = "Démarrage de l'app"
$mobile.Range("B4").Value = "cliquer l'icon"
$mobile.Range("C4").Value = "splash screen`napparition du login"
$mobile.Range("D4").Value = "ok"

# --- Row 10 : login (filled early by the author, before rows 5-9) ---
$mobile.Range("A10").Value = "login"

# --- Row 5 : conditions generales ---
$mobile.Range("A5").Value = "conditions generales"
$mobile.Range("B5").Value = "cliquer le bouton conditions generlaes"
$mobile.Range("C5").Value = "chargement`napparition des conditions generales"
$mobile.Range("D5").Value = "ok"

# --- Row 6 : multi-lingue ---
$mobile.Range("A6").Value = "multi-lingue"
$mobile.Range("B6").Value = "cliquer l'icon retour`nchoisir la langue malgache`ncliquer le bouton conditions generlaes"
$mobile.Range("C6").Value = "l'app est en malgache`nchargement`napparition des conditions generales en malgaches"
$mobile.Range("D6").Value = "ok"

# --- Row 10 (continued) ---
$mobile.Range("B10").Value = "remplir le formulaire par`nemail: teste@gmail.com`nmots de passe: testestes"
$mobile.Range("C10").Value = "chargement de la fenêtre`nla fenêtre de bienvenue apparait"

# --- Row 7 : inscription (highlighted, filled last) ---
$mobile.Range("A7").Value = "inscription"
$mobile.Range("B7").Value = "cliquer creer un compte`nremplir le formulaire par`nnom complet: `nemail: teste@gmail.com`nmots de passe: testestes"

# --- Row 2 : author ---
$mobile.Range("A2").Value = "Auteur"
$mobile.Range("B2").Value = "Anah"

# --- Row 3 : column headers ---
$mobile.Range("A3").Value = "Nom de la fonctionnalité"
$mobile.Range("B3").Value = "Étapes d’exécution"
$mobile.Range("C3").Value = "Résultat attendu"
$mobile.Range("D3").Value = "Résultat réel"

# --- formatting pass ---
$mobile.Range("A1").Font.Bold = $true
$mobile.Range("B1:D1").Font.Bold = $false

$mobile.Range("A2").Font.Bold = $true
$mobile.Range("B2:D2").Font.Bold = $false

$mobile.Range("A3").Font.Bold = $true
$mobile.Range("B3").Font.Bold = $true
$mobile.Range("B3").WrapText = $true
$mobile.Range("C3").Font.Bold = $true
$mobile.Range("D3").Font.Bold = $true
$mobile.Rows.Item(3).RowHeight = 45

$mobile.Range("A4").Font.Bold = $true
$mobile.Range("C4").WrapText = $true
$mobile.Rows.Item(4).RowHeight = 30

$mobile.Range("A5").Font.Bold = $true
$mobile.Range("C5").WrapText = $true
$mobile.Rows.Item(5).RowHeight = 30

$mobile.Range("A6").Font.Bold = $true
$mobile.Range("B6").WrapText = $true
$mobile.Range("C6").WrapText = $true
$mobile.Rows.Item(6).RowHeight = 45

$mobile.Range("A7").Font.Bold = $true
$mobile.Range("A7").Interior.Color = 65535
$mobile.Range("B7").WrapText = $true
$mobile.Rows.Item(7).RowHeight = 75

$mobile.Range("A10").Font.Bold = $true
$mobile.Range("B10").WrapText = $true
$mobile.Range("C10").WrapText = $true
$mobile.Rows.Item(10).RowHeight = 45

# --- activate "mobile" and select A7, matching the saved view state ---
$mobile.Activate()
$mobile.Range("A7").Select()
